$wb = $excel.ActiveWorkbook

# Sheet 1: "Full results"
$ws1 = $wb.Worksheets.Item("Full results")

# Row 2 - NULL MODEL
$ws1.Range("C2").Value = 0.908941580497606
$ws1.Range("D2").Value = 0.0911546741138437
$ws1.Range("E2").Value = 1.00009625461145
$ws1.Range("J2").Value = 0.0911459009005673
$ws1.Range("K2").Value = 0.0943672761142998
$ws1.Range("L2").Value = 0.00261650202457523
$ws1.Range("M2").Value = 0.0370888117502591
$ws1.Range("N2").Value = 0.0969837781388751

# Row 3 - CONDITIONAL MODEL
$ws1.Range("F3").Value = 0.874465952653143
$ws1.Range("G3").Value = 0.0943763593997957

# Row 4 - COMPLETE MODEL
$ws1.Range("H4").Value = 0.871849198778182
$ws1.Range("I4").Value = 0.0876149474324018
$ws1.Range("O4").Value = 0.128234712650826

# Sheet 2: "For plotting"
$ws2 = $wb.Worksheets.Item("For plotting")

# Row 2 - Sibcorr
$ws2.Range("C2").Value = 0.0911459009005673
$ws2.Range("D2").Value = 0.0360600100989626
$ws2.Range("E2").Value = 0.146231791702172
$ws2.Range("F2").Value = 948

# Row 3 - IOLIB
$ws2.Range("C3").Value = 0.0969837781388751
$ws2.Range("D3").Value = 0.0212652213961603
$ws2.Range("E3").Value = 0.17270233488159
$ws2.Range("F3").Value = 948

# Row 4 - IORAD
$ws2.Range("C4").Value = 0.128234712650826
$ws2.Range("D4").Value = 0.0381265753423746
$ws2.Range("E4").Value = 0.218342849959278
$ws2.Range("F4").Value = 948
